# Actualizacion Datos Personales 4 nov
#
# This script fills in the previously-missing contact details for two
# students whose rows only had their identification columns (A-D) populated.
#
# Sheet "3AEM" (sheet1), row 7  -> CRISTIAN JAVIER CORTEZ ANTONIO
# Sheet "3APM" (sheet3), row 6  -> CESAR CUEVAS CUATRA

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 3AEM - row 7
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("3AEM")

$ws1.Range("E7").Value = "cristianantoniof2020@gmail.com"

$ws1.Range("F7").Value = "'2721079631"
$ws1.Range("F7").Style = "Normal"

$ws1.Range("G7").Value = "'2721079631"
$ws1.Range("G7").Style = "Normal"

$ws1.Range("H7").Value = "FRANCISCO JAVIER CORTÉZ LEYNES"

$ws1.Range("I7").Value = "Cristianantoniof2020@gmail.com"

# ---------------------------------------------------------------------------
# Sheet 3APM - row 6
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("3APM")

$ws3.Range("E6").Value = "cesarcuevasc3@gmail.com"

$ws3.Range("F6").Value = "'2722848082"
$ws3.Range("F6").Style = "Normal"

$ws3.Range("G6").Value = "'2722848082"
$ws3.Range("G6").Style = "Normal"

$ws3.Range("H6").Value = "GERARDO CUEVAS MACUIXTLE"

$ws3.Range("I6").Value = "cesarcuevasc3@gmail.com"

$ws3.Range("J6").Value = "'2722848082"
$ws3.Range("J6").Style = "Normal"
